$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A36").Value = "HGP1Q2"
$ws.Range("B36").Value = "Chip Epson"
$ws.Range("C36").Value = "C9345"
$ws.Range("D36").Value = 0
$ws.Range("E36").Value = 100000
$ws.Range("F36").Value = 20
$ws.Range("G36").Value = 0
$ws.Range("H36").Formula = "=(E36-D36)*G36"
$ws.Range("I36").Formula = "=D36*F36"
$ws.Range("J36").Value = 0
